$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the per-row data that gets reshuffled between rows 2..12
# A=1 (Id), B=2 (Taxonsorteringsordning), D=4 (Rodlistade), E=5 (TaxonId),
# F=6 (Artnamn), G=7 (Vetenskapligt namn), H=8 (Auktor),
# Q=17 (Ost), R=18 (Nord), S=19 (Noggrannhet)
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18, 19)

# Mapping: new row N gets the values that currently live in row Map[N] (before state)
$map = @{
    2  = 6
    3  = 2
    4  = 3
    5  = 4
    6  = 12
    7  = 11
    8  = 8
    9  = 10
    10 = 9
    11 = 5
    12 = 7
}

# Snapshot current ("before") values for every row/column involved so that
# writes don't clobber data we still need to read later.
$snapshot = @{}
foreach ($r in 2..12) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Apply the permutation.
foreach ($r in 2..12) {
    $srcRow = $map[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
